$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1749153331348/OwQjboHWXD-camera-video-segment-driver-1749153333848.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPCYFEHDQ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T200219Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCcsOvH6sVCM7zznQgEqMN5Du%2F%2FOWeev7uYddPJZj7RtAIgUzq2l3zAXMsS9QBz%2BfjNVvRkK%2BPVT2hpSq9adP5FNkkq3QMIYxAEGgw3ODEyMDQ5NDIyNDQiDNSYBD651hX2ub8S8Sq6A6OEQfMHm20%2F8ALVib9rHcOii7vCtzmshBCR%2Ftk5jvmxdoe0BblyUbv5vuHApX36Bzj3nG7E9nRNEBKiKH4z7yiTOSOVNG2ltKGLS1ZoJ3rVDjxci92G%2FgE8scsAARw2QicmOBk0rel9vZsGVArasuXyhNUp%2FFtrc3dmKbDyyKQQ%2BgJXFiCG9Rwfj8KJBMgMWc8fQJfUqBztR80MnHDLc97sWaW3afZ9S3Q9U3IidjF7TPbFUO2KWUZWbeYSzfcU8Gib5JiKLhiefYthgXp2MGjNjZt16HtzMPkhrb1Ht8GNvnYNDHyVGB%2B4JCDEQ2eOSAqLGkwH0Oig6h6Ztm0cj48YfgLY6i9jxoCKU6%2FPLrbkFRX4FtFcwJCaJfoU%2FSu8u0kEO%2FAWzKfsy2gUeEwZbx6zX39umuTbfJ2kl%2FTamOOL9x2z2TPnOlSGVSqL4XUcBMOun%2B3t9hKXmmfKvk%2FXfDRDI0G2roy%2FR09zVrgAM0qLHvlBHVeZUs8NrajhTP9lBu2ZRrutHfbDiyTnesm3MGm8r28ji1M6kLbBlGhDFsrxyfTsaUnJNr6anCver9AOEI79GWkmmtL5kqowqe2MwgY6pQFO7O9XVd9e83oUDScEy19DJ13f8Bd3pCXc47XkvKAk5kSpV7CVif4G%2FtM1kUHdtDjI2J5aWwmkzjcxfqgva8tCn4sjXSstESfTf8ERA8WIkjEczPaUbgNmsyoWo96Mlc2vduzaTbGVKBi%2B%2FJI4oAFfPKXYL5ZWvqnbAGVD5Li6%2FFgdGFID0heaIEqrhuC9fQ88h1xKwD5EBe71DLbR6eRnvdDUxFE%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2004%3A02%3A19%20GMT&X-Amz-Signature=2c9b3cb661b3d04dbf62ae65ff532ae39c5f2ba9729419c3e8227fc9041558ea'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1749153331348/oV5nVwJUOa-camera-video-segment-1749153333848.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPCYFEHDQ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T200219Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCcsOvH6sVCM7zznQgEqMN5Du%2F%2FOWeev7uYddPJZj7RtAIgUzq2l3zAXMsS9QBz%2BfjNVvRkK%2BPVT2hpSq9adP5FNkkq3QMIYxAEGgw3ODEyMDQ5NDIyNDQiDNSYBD651hX2ub8S8Sq6A6OEQfMHm20%2F8ALVib9rHcOii7vCtzmshBCR%2Ftk5jvmxdoe0BblyUbv5vuHApX36Bzj3nG7E9nRNEBKiKH4z7yiTOSOVNG2ltKGLS1ZoJ3rVDjxci92G%2FgE8scsAARw2QicmOBk0rel9vZsGVArasuXyhNUp%2FFtrc3dmKbDyyKQQ%2BgJXFiCG9Rwfj8KJBMgMWc8fQJfUqBztR80MnHDLc97sWaW3afZ9S3Q9U3IidjF7TPbFUO2KWUZWbeYSzfcU8Gib5JiKLhiefYthgXp2MGjNjZt16HtzMPkhrb1Ht8GNvnYNDHyVGB%2B4JCDEQ2eOSAqLGkwH0Oig6h6Ztm0cj48YfgLY6i9jxoCKU6%2FPLrbkFRX4FtFcwJCaJfoU%2FSu8u0kEO%2FAWzKfsy2gUeEwZbx6zX39umuTbfJ2kl%2FTamOOL9x2z2TPnOlSGVSqL4XUcBMOun%2B3t9hKXmmfKvk%2FXfDRDI0G2roy%2FR09zVrgAM0qLHvlBHVeZUs8NrajhTP9lBu2ZRrutHfbDiyTnesm3MGm8r28ji1M6kLbBlGhDFsrxyfTsaUnJNr6anCver9AOEI79GWkmmtL5kqowqe2MwgY6pQFO7O9XVd9e83oUDScEy19DJ13f8Bd3pCXc47XkvKAk5kSpV7CVif4G%2FtM1kUHdtDjI2J5aWwmkzjcxfqgva8tCn4sjXSstESfTf8ERA8WIkjEczPaUbgNmsyoWo96Mlc2vduzaTbGVKBi%2B%2FJI4oAFfPKXYL5ZWvqnbAGVD5Li6%2FFgdGFID0heaIEqrhuC9fQ88h1xKwD5EBe71DLbR6eRnvdDUxFE%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2004%3A02%3A19%20GMT&X-Amz-Signature=f440fd04ef12ab5caf59e37a82b6160c944f8b3822c377487a6be71f0b1a41e4'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1749147387468/GKgOu1ZeR3-camera-video-segment-driver-1749147389968.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPCYFEHDQ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T200219Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCcsOvH6sVCM7zznQgEqMN5Du%2F%2FOWeev7uYddPJZj7RtAIgUzq2l3zAXMsS9QBz%2BfjNVvRkK%2BPVT2hpSq9adP5FNkkq3QMIYxAEGgw3ODEyMDQ5NDIyNDQiDNSYBD651hX2ub8S8Sq6A6OEQfMHm20%2F8ALVib9rHcOii7vCtzmshBCR%2Ftk5jvmxdoe0BblyUbv5vuHApX36Bzj3nG7E9nRNEBKiKH4z7yiTOSOVNG2ltKGLS1ZoJ3rVDjxci92G%2FgE8scsAARw2QicmOBk0rel9vZsGVArasuXyhNUp%2FFtrc3dmKbDyyKQQ%2BgJXFiCG9Rwfj8KJBMgMWc8fQJfUqBztR80MnHDLc97sWaW3afZ9S3Q9U3IidjF7TPbFUO2KWUZWbeYSzfcU8Gib5JiKLhiefYthgXp2MGjNjZt16HtzMPkhrb1Ht8GNvnYNDHyVGB%2B4JCDEQ2eOSAqLGkwH0Oig6h6Ztm0cj48YfgLY6i9jxoCKU6%2FPLrbkFRX4FtFcwJCaJfoU%2FSu8u0kEO%2FAWzKfsy2gUeEwZbx6zX39umuTbfJ2kl%2FTamOOL9x2z2TPnOlSGVSqL4XUcBMOun%2B3t9hKXmmfKvk%2FXfDRDI0G2roy%2FR09zVrgAM0qLHvlBHVeZUs8NrajhTP9lBu2ZRrutHfbDiyTnesm3MGm8r28ji1M6kLbBlGhDFsrxyfTsaUnJNr6anCver9AOEI79GWkmmtL5kqowqe2MwgY6pQFO7O9XVd9e83oUDScEy19DJ13f8Bd3pCXc47XkvKAk5kSpV7CVif4G%2FtM1kUHdtDjI2J5aWwmkzjcxfqgva8tCn4sjXSstESfTf8ERA8WIkjEczPaUbgNmsyoWo96Mlc2vduzaTbGVKBi%2B%2FJI4oAFfPKXYL5ZWvqnbAGVD5Li6%2FFgdGFID0heaIEqrhuC9fQ88h1xKwD5EBe71DLbR6eRnvdDUxFE%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2004%3A02%3A19%20GMT&X-Amz-Signature=71bc8fb6d501637a24908542744c849d9e7b38c693945c9979684fbd6d69214c'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1749147387468/dQoeD1EUOr-camera-video-segment-1749147389968.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPCYFEHDQ%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T200219Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCcsOvH6sVCM7zznQgEqMN5Du%2F%2FOWeev7uYddPJZj7RtAIgUzq2l3zAXMsS9QBz%2BfjNVvRkK%2BPVT2hpSq9adP5FNkkq3QMIYxAEGgw3ODEyMDQ5NDIyNDQiDNSYBD651hX2ub8S8Sq6A6OEQfMHm20%2F8ALVib9rHcOii7vCtzmshBCR%2Ftk5jvmxdoe0BblyUbv5vuHApX36Bzj3nG7E9nRNEBKiKH4z7yiTOSOVNG2ltKGLS1ZoJ3rVDjxci92G%2FgE8scsAARw2QicmOBk0rel9vZsGVArasuXyhNUp%2FFtrc3dmKbDyyKQQ%2BgJXFiCG9Rwfj8KJBMgMWc8fQJfUqBztR80MnHDLc97sWaW3afZ9S3Q9U3IidjF7TPbFUO2KWUZWbeYSzfcU8Gib5JiKLhiefYthgXp2MGjNjZt16HtzMPkhrb1Ht8GNvnYNDHyVGB%2B4JCDEQ2eOSAqLGkwH0Oig6h6Ztm0cj48YfgLY6i9jxoCKU6%2FPLrbkFRX4FtFcwJCaJfoU%2FSu8u0kEO%2FAWzKfsy2gUeEwZbx6zX39umuTbfJ2kl%2FTamOOL9x2z2TPnOlSGVSqL4XUcBMOun%2B3t9hKXmmfKvk%2FXfDRDI0G2roy%2FR09zVrgAM0qLHvlBHVeZUs8NrajhTP9lBu2ZRrutHfbDiyTnesm3MGm8r28ji1M6kLbBlGhDFsrxyfTsaUnJNr6anCver9AOEI79GWkmmtL5kqowqe2MwgY6pQFO7O9XVd9e83oUDScEy19DJ13f8Bd3pCXc47XkvKAk5kSpV7CVif4G%2FtM1kUHdtDjI2J5aWwmkzjcxfqgva8tCn4sjXSstESfTf8ERA8WIkjEczPaUbgNmsyoWo96Mlc2vduzaTbGVKBi%2B%2FJI4oAFfPKXYL5ZWvqnbAGVD5Li6%2FFgdGFID0heaIEqrhuC9fQ88h1xKwD5EBe71DLbR6eRnvdDUxFE%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2004%3A02%3A19%20GMT&X-Amz-Signature=83dccfa7c910ba6298ebcb62cfa74a1ed29e6fc38f8c409138cc01f2acf7345b'
